$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 38087.5
$ws.Range("J3").Value = 38087.5
$ws.Range("L3").Value = 38087.5
$ws.Range("N3").Value = -38315.5

$ws.Range("H40").Value = 11266
$ws.Range("I40").Value = 15149
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 15149
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -14974
$ws.Range("N40").Value = -3850

$ws.Range("H93").Value = 34051.414
$ws.Range("J93").Value = 34051.414
$ws.Range("L93").Value = 34051.414
$ws.Range("N93").Value = -39043.414

$ws.Range("H95").Value = 33318
$ws.Range("J95").Value = 33318
$ws.Range("L95").Value = 33318
$ws.Range("N95").Value = -38810

$ws.Range("H102").Value = 38087.5
$ws.Range("J102").Value = 38087.5
$ws.Range("L102").Value = 38087.5
$ws.Range("N102").Value = -44577.5

$ws.Range("H105").Value = 38664
$ws.Range("J105").Value = 38664
$ws.Range("L105").Value = 38664
$ws.Range("N105").Value = -45652

$ws.Range("H123").Value = 39998.184
$ws.Range("J123").Value = 39998.184
$ws.Range("L123").Value = 39998.184
$ws.Range("N123").Value = -49798.184

$ws.Range("H129").Value = 4963.467
$ws.Range("I129").Value = 4174.25
$ws.Range("J129").Value = 5250.4546
$ws.Range("K129").Value = 12522.75
$ws.Range("L129").Value = 15751.3638
$ws.Range("M129").Value = -7522.75
$ws.Range("N129").Value = -25751.3638

$ws.Range("H132").Value = 19682.21
$ws.Range("I132").Value = 3122.5107
$ws.Range("J132").Value = 175343.4
$ws.Range("K132").Value = 9367.5321
$ws.Range("L132").Value = 526030.2
$ws.Range("M132").Value = -6837.5321
$ws.Range("N132").Value = -531090.2

$ws.Range("H135").Value = 18519690
$ws.Range("I135").Value = 1228.2727
$ws.Range("J135").Value = 100000920
$ws.Range("K135").Value = 11054.4543
$ws.Range("L135").Value = 900008280
$ws.Range("M135").Value = -8519.454299999999
$ws.Range("N135").Value = -900013350

$ws.Range("H138").Value = 2913.4727
$ws.Range("J138").Value = 3204.5938
$ws.Range("L138").Value = 9613.7814
$ws.Range("N138").Value = -19893.7814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 35000
$ws.Range("J24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("N24").Value = -35748

$ws.Range("H25").Value = 1696.8334
$ws.Range("I25").Value = 436.2
$ws.Range("J25").Value = 8000
$ws.Range("K25").Value = 436.2
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = -34.19999999999999
$ws.Range("N25").Value = -8804

$ws.Range("H32").Value = 29926.984
$ws.Range("I32").Value = 31394.492
$ws.Range("J32").Value = 15496.5
$ws.Range("K32").Value = 31394.492
$ws.Range("L32").Value = 15496.5
$ws.Range("M32").Value = -31107.492
$ws.Range("N32").Value = -16070.5

$ws.Range("H95").Value = 27600
$ws.Range("J95").Value = 27600
$ws.Range("L95").Value = 27600
$ws.Range("N95").Value = -33092

$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164

$ws.Range("H103").Value = 34358
$ws.Range("J103").Value = 34358
$ws.Range("L103").Value = 34358
$ws.Range("N103").Value = -36702

$ws.Range("H105").Value = 49181
$ws.Range("J105").Value = 49181
$ws.Range("L105").Value = 49181
$ws.Range("N105").Value = -56169

$ws.Range("H106").Value = 48367.332
$ws.Range("J106").Value = 48367.332
$ws.Range("L106").Value = 48367.332
$ws.Range("N106").Value = -50891.332

$ws.Range("H109").Value = 41877
$ws.Range("J109").Value = 41877
$ws.Range("L109").Value = 41877
$ws.Range("N109").Value = -44651

$ws.Range("H114").Value = 34145.5
$ws.Range("J114").Value = 34145.5
$ws.Range("L114").Value = 34145.5
$ws.Range("N114").Value = -42823.5

$ws.Range("H122").Value = 1716.5
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H130").Value = 43747.332
$ws.Range("J130").Value = 43747.332
$ws.Range("L130").Value = 43747.332
$ws.Range("N130").Value = -53787.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 44724.668
$ws.Range("J95").Value = 44724.668
$ws.Range("L95").Value = 44724.668
$ws.Range("N95").Value = -50216.668

$ws.Range("H119").Value = 46711
$ws.Range("J119").Value = 46711
$ws.Range("L119").Value = 46711
$ws.Range("N119").Value = -56387

$ws.Range("H130").Value = 48367.332
$ws.Range("J130").Value = 48367.332
$ws.Range("L130").Value = 48367.332
$ws.Range("N130").Value = -58407.332

$ws.Range("H134").Value = 3543.2273
$ws.Range("I134").Value = 3246.75
$ws.Range("J134").Value = 3899
$ws.Range("K134").Value = 9740.25
$ws.Range("L134").Value = 11697
$ws.Range("M134").Value = -7205.25
$ws.Range("N134").Value = -16767

$ws.Range("H135").Value = 29736.389
$ws.Range("J135").Value = 29736.389
$ws.Range("L135").Value = 29736.389
$ws.Range("N135").Value = -39876.389

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8650.92
$ws.Range("I31").Value = 4780.143
$ws.Range("K31").Value = 4780.143
$ws.Range("M31").Value = -4485.143

$ws.Range("H34").Value = 8650.92
$ws.Range("I34").Value = 4780.143
$ws.Range("K34").Value = 4780.143
$ws.Range("M34").Value = -4578.143

$ws.Range("H43").Value = 28216.334
$ws.Range("J43").Value = 28216.334
$ws.Range("L43").Value = 28216.334
$ws.Range("N43").Value = -28584.334

$ws.Range("H101").Value = 28216.334
$ws.Range("J101").Value = 28216.334
$ws.Range("L101").Value = 28216.334
$ws.Range("N101").Value = -34706.334

$ws.Range("H106").Value = 45980
$ws.Range("J106").Value = 45980
$ws.Range("L106").Value = 45980
$ws.Range("N106").Value = -48504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3900
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8800

$ws.Range("H123").Value = 15660.3
$ws.Range("J123").Value = 15660.3
$ws.Range("L123").Value = 15660.3
$ws.Range("N123").Value = -20560.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1981.6
$ws.Range("I46").Value = 1236.1111
$ws.Range("J46").Value = 3099.8333
$ws.Range("K46").Value = 1236.1111
$ws.Range("L46").Value = 3099.8333
$ws.Range("M46").Value = -1048.1111
$ws.Range("N46").Value = -3475.8333

$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 944.44446
$ws.Range("J55").Value = 1125
$ws.Range("K55").Value = 944.44446
$ws.Range("L55").Value = 1125
$ws.Range("M55").Value = -771.44446
$ws.Range("N55").Value = -1471

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H68").Value = 2395.6875
$ws.Range("I68").Value = 2310.077
$ws.Range("K68").Value = 2310.077
$ws.Range("M68").Value = -1561.077

$ws.Range("H71").Value = 2395.6875
$ws.Range("I71").Value = 2310.077
$ws.Range("K71").Value = 11550.385
$ws.Range("M71").Value = -7806.385000000002

$ws.Range("H94").Value = 54650
$ws.Range("J94").Value = 54650
$ws.Range("L94").Value = 54650
$ws.Range("N94").Value = -56002

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H121").Value = 44406.8
$ws.Range("J121").Value = 44406.8
$ws.Range("L121").Value = 44406.8
$ws.Range("N121").Value = -47900.8

$ws.Range("H127").Value = 50715
$ws.Range("J127").Value = 50715
$ws.Range("L127").Value = 50715
$ws.Range("N127").Value = -60635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 97887.75
$ws.Range("J92").Value = 97887.75
$ws.Range("L92").Value = 97887.75
$ws.Range("N92").Value = -102879.75

$ws.Range("H97").Value = 38186
$ws.Range("J97").Value = 38186
$ws.Range("L97").Value = 38186
$ws.Range("N97").Value = -40168

$ws.Range("H103").Value = 41814.668
$ws.Range("J103").Value = 41814.668
$ws.Range("L103").Value = 41814.668
$ws.Range("N103").Value = -44158.668

$ws.Range("H104").Value = 43681
$ws.Range("J104").Value = 43681
$ws.Range("L104").Value = 43681
$ws.Range("N104").Value = -50669

$ws.Range("H137").Value = 67660
$ws.Range("J137").Value = 67660
$ws.Range("L137").Value = 67660
$ws.Range("N137").Value = -77860

